$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")

# The new value "1" looks numeric, but the original column stores rule
# names as text (e.g. "R10".."R40"), so it must stay a text cell (shared
# string) rather than turn into a number. Typing it directly would make
# Excel store it as a real number, and prefixing with an apostrophe would
# flip on the "Quote Prefix" cell flag (changing the cell style). Instead,
# stage the text value in a scratch cell that is already blank and inside
# the used range (so dimension/row spans are not perturbed), copy it, and
# paste only the value into B11 - this preserves B11's existing style
# (s="23") while making the stored value a genuine text string.
$helper = $ws.Range("B5")
$helper.Formula = "=""1"""
$helper.Copy()
$target.PasteSpecial(-4163)   # xlPasteValues
$helper.ClearContents()

$wb.Save()
